$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 76.25
$ws.Range("I4").Value = 76.25
$ws.Range("K4").Value = 76.25
$ws.Range("M4").Value = 37.75

$ws.Range("H6").Value = 650
$ws.Range("I6").Value = 1000
$ws.Range("K6").Value = 3000
$ws.Range("M6").Value = -2888

$ws.Range("H11").Value = 183.85185
$ws.Range("I11").Value = 183.85185
$ws.Range("K11").Value = 183.85185
$ws.Range("M11").Value = -43.85185000000001

$ws.Range("H17").Value = 2214
$ws.Range("I17").Value = 2361
$ws.Range("J17").Value = 2195.625
$ws.Range("K17").Value = 7083
$ws.Range("L17").Value = 6586.875
$ws.Range("M17").Value = -6915
$ws.Range("N17").Value = -6922.875

$ws.Range("H19").Value = 688
$ws.Range("J19").Value = 687
$ws.Range("L19").Value = 687
$ws.Range("N19").Value = -1037

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = ""

$ws.Range("H26").Value = 2000
$ws.Range("J26").Value = 2000
$ws.Range("L26").Value = 2000
$ws.Range("N26").Value = -2688

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = ""

$ws.Range("H43").Value = 3666.6667
$ws.Range("J43").Value = 3400
$ws.Range("L43").Value = 3400
$ws.Range("N43").Value = -3538

$ws.Range("H112").Value = 1712.2258
$ws.Range("J112").Value = 1727.1333
$ws.Range("L112").Value = 5181.3999
$ws.Range("N112").Value = -7397.3999

$ws.Range("H116").Value = 9029.652
$ws.Range("I116").Value = 8271
$ws.Range("J116").Value = 9517.357
$ws.Range("K116").Value = 8271
$ws.Range("L116").Value = 9517.357
$ws.Range("M116").Value = -4829
$ws.Range("N116").Value = -16401.357

$ws.Range("H130").Value = 69997.27
$ws.Range("J130").Value = 69997.27
$ws.Range("L130").Value = 69997.27
$ws.Range("N130").Value = -80037.27

$ws.Range("H138").Value = 3437.28
$ws.Range("I138").Value = 2470.389
$ws.Range("J138").Value = 3981.1562
$ws.Range("K138").Value = 7411.167
$ws.Range("L138").Value = 11943.4686
$ws.Range("M138").Value = -2271.167
$ws.Range("N138").Value = -22223.4686

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3437.4614
$ws.Range("I32").Value = 3539.6458
$ws.Range("J32").Value = 2211.25
$ws.Range("K32").Value = 3539.6458
$ws.Range("L32").Value = 2211.25
$ws.Range("M32").Value = -3252.6458
$ws.Range("N32").Value = -2785.25

$ws.Range("H74").Value = 1204.3721
$ws.Range("I74").Value = 1145.8462
$ws.Range("J74").Value = 1775
$ws.Range("K74").Value = 1145.8462
$ws.Range("L74").Value = 1775
$ws.Range("M74").Value = -271.8462
$ws.Range("N74").Value = -3523

$ws.Range("H77").Value = 1204.3721
$ws.Range("I77").Value = 1145.8462
$ws.Range("J77").Value = 1775
$ws.Range("K77").Value = 5729.231
$ws.Range("L77").Value = 8875
$ws.Range("M77").Value = -1361.231
$ws.Range("N77").Value = -17611

$ws.Range("H132").Value = 2663.3333
$ws.Range("I132").Value = 2337.2727
$ws.Range("K132").Value = 7011.8181
$ws.Range("M132").Value = -4481.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1847.0834
$ws.Range("I105").Value = 1640.5555
$ws.Range("K105").Value = 1640.5555
$ws.Range("M105").Value = 106.4445000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2155.9167
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = ""

$ws.Range("H34").Value = 2155.9167
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = ""

$ws.Range("H38").Value = 11935.4
$ws.Range("I38").Value = 11935.4
$ws.Range("K38").Value = 11935.4
$ws.Range("M38").Value = -11558.4

$ws.Range("H46").Value = 11935.4
$ws.Range("I46").Value = 11935.4
$ws.Range("K46").Value = 11935.4
$ws.Range("M46").Value = -11724.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 612.8182
$ws.Range("I5").Value = 393.58823
$ws.Range("J5").Value = 1358.2
$ws.Range("K5").Value = 1180.76469
$ws.Range("L5").Value = 4074.6
$ws.Range("M5").Value = -1068.76469
$ws.Range("N5").Value = -4298.6

$ws.Range("H14").Value = 41
$ws.Range("I14").Value = 41
$ws.Range("K14").Value = 123
$ws.Range("M14").Value = 50

$ws.Range("H113").Value = 1836.375
$ws.Range("I113").Value = 2932.3333
$ws.Range("J113").Value = 1178.8
$ws.Range("K113").Value = 8796.999899999999
$ws.Range("L113").Value = 3536.4
$ws.Range("M113").Value = -6626.999899999999
$ws.Range("N113").Value = -7876.4

$ws.Range("H122").Value = 8582.091
$ws.Range("J122").Value = 14999.833
$ws.Range("L122").Value = 134998.497
$ws.Range("N122").Value = -139898.497

$ws.Range("H131").Value = 2660.3845
$ws.Range("I131").Value = 2314.25
$ws.Range("J131").Value = 3214.2
$ws.Range("K131").Value = 6942.75
$ws.Range("L131").Value = 9642.599999999999
$ws.Range("M131").Value = -1902.75
$ws.Range("N131").Value = -19722.6

$ws.Range("H135").Value = 612.8182
$ws.Range("I135").Value = 393.58823
$ws.Range("J135").Value = 1358.2
$ws.Range("K135").Value = 3542.29407
$ws.Range("L135").Value = 12223.8
$ws.Range("M135").Value = -1007.29407
$ws.Range("N135").Value = -17293.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 62036.3
$ws.Range("I80").Value = 106093.63
$ws.Range("K80").Value = 106093.63
$ws.Range("M80").Value = -105095.63

$ws.Range("H83").Value = 62036.3
$ws.Range("I83").Value = 106093.63
$ws.Range("K83").Value = 530468.15
$ws.Range("M83").Value = -525476.15

$ws.Range("H112").Value = 48500
$ws.Range("J112").Value = 48500
$ws.Range("L112").Value = 48500
$ws.Range("N112").Value = -50716

$ws.Range("H126").Value = 5307.8
$ws.Range("I126").Value = 5341.6665
$ws.Range("K126").Value = 16024.9995
$ws.Range("M126").Value = -13554.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 911.4286
$ws.Range("I55").Value = 897
$ws.Range("J55").Value = 998
$ws.Range("K55").Value = 897
$ws.Range("L55").Value = 998
$ws.Range("M55").Value = -724
$ws.Range("N55").Value = -1344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3170
$ws.Range("I62").Value = 2125
$ws.Range("K62").Value = 2125
$ws.Range("M62").Value = -1501

$ws.Range("H65").Value = 3170
$ws.Range("I65").Value = 2125
$ws.Range("K65").Value = 10625
$ws.Range("M65").Value = -7505

$ws.Range("H136").Value = 1341
$ws.Range("I136").Value = 1306.1333
$ws.Range("K136").Value = 3918.3999
$ws.Range("M136").Value = -1368.3999
